$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

$wsTraining = $wb.Worksheets.Item("Training")
$wsValidation = $wb.Worksheets.Item("Validation")
$wsExperimentation = $wb.Worksheets.Item("Experimentation")
$wsAdjusted = $wb.Worksheets.Item("Adjusted")

# ---------------------------------------------------------------------------
# 1. Adjusted sheet: the "Multilayer Perceptron (500 epochs / 51 hidden
#    nodes)" placeholder row becomes a "Multilayer Perceptron" row with
#    its own metrics, and it no longer needs the taller, wrapped row
#    height. (Label written first so the shared-string table allocates
#    "Multilayer Perceptron" before the longer training/validation label.)
# ---------------------------------------------------------------------------
$wsAdjusted.Cells.Item(4,1).Value = "Multilayer Perceptron"

# ---------------------------------------------------------------------------
# 2. Training sheet: insert a new "Multilayer Perceptron (10 epochs / 537
#    hidden nodes)" row above the existing 500-epoch row, and fill in its
#    metrics.
# ---------------------------------------------------------------------------
$wsTraining.Rows.Item(7).Insert()

$wsTraining.Cells.Item(7,1).Value = "Multilayer Perceptron (10 epochs / 537 hidden nodes)"
$wsTraining.Cells.Item(7,2).Value = 0.884
$wsTraining.Cells.Item(7,2).NumberFormat = "0.000"
$wsTraining.Cells.Item(7,3).Value = 0.893
$wsTraining.Cells.Item(7,4).Value = 0.884
$wsTraining.Cells.Item(7,5).Value = 0.885
$wsTraining.Rows.Item(7).RowHeight = 29

$wsTraining.Range("B7:E7").Select()

# ---------------------------------------------------------------------------
# 3. Validation sheet: the "Multilayer Perceptron (25 epochs / 537 hidden
#    nodes)" label row had no metrics yet - rename it to the 10-epoch
#    variant (re-using the shared string created in step 2) and fill in
#    its metrics.
# ---------------------------------------------------------------------------
$wsValidation.Cells.Item(6,1).Value = "Multilayer Perceptron (10 epochs / 537 hidden nodes)"
$wsValidation.Cells.Item(6,2).Value = 0.595
$wsValidation.Cells.Item(6,3).Value = 0.677
$wsValidation.Cells.Item(6,4).Value = 0.595
$wsValidation.Cells.Item(6,5).Value = 0.584

$wsValidation.Range("B3:E3").Select()

# ---------------------------------------------------------------------------
# 4. Experimentation sheet: scroll position tweak, it is no longer the
#    active tab (Adjusted becomes active instead - handled below).
# ---------------------------------------------------------------------------
$wsExperimentation.Range("A30").Select()

# ---------------------------------------------------------------------------
# 5. Back to the Adjusted sheet: finish filling in the new row's metrics
#    and note (the "decay = true" string is allocated last).
# ---------------------------------------------------------------------------
$wsAdjusted.Cells.Item(4,2).Value = 0.646
$wsAdjusted.Cells.Item(4,2).NumberFormat = "0.000"
$wsAdjusted.Cells.Item(4,3).Value = 0.682
$wsAdjusted.Cells.Item(4,4).Value = 0.646
$wsAdjusted.Cells.Item(4,5).Value = 0.644
$wsAdjusted.Cells.Item(4,6).Value = "decay = true"
$wsAdjusted.Rows.Item(4).AutoFit()

# ---------------------------------------------------------------------------
# 5. Drop the "Testing" sheet entirely (it was just a blank template) -
#    deleting it makes "Adjusted" (the new last sheet) active, matching the
#    desired tab selection.
# ---------------------------------------------------------------------------
[void]$wb.Worksheets.Item("Testing").Delete()
$wsAdjusted.Activate()
